$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), formatted like the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# New "Save" value for the data row (H2)
$ws.Range("H2").Value = 1
